# Refresh the crypto price/volume snapshot (Price column D, Volume(1h) column E).
# Values that look like plain numbers are entered with a leading apostrophe so
# Excel keeps them as text (matching the sheet's existing inline-string cells)
# instead of silently re-parsing "11.30" -> 11.3, "0.3376" -> 0.3376 as a Double.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "22.480.75"
$ws.Range("D3").Value = "1.571.67"
$ws.Range("E3").Value = "  +0.27%  "
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("E5").Value = "  +0.03%  "
$ws.Range("D6").Value = "'290.99"
$ws.Range("E6").Value = "  -0.03%  "
$ws.Range("E7").Value = "  -1.40%  "
$ws.Range("D8").Value = "'49.97"
$ws.Range("E8").Value = "  +2.06%  "
$ws.Range("D9").Value = "'0.3376"
$ws.Range("E9").Value = "  -0.59%  "
$ws.Range("D10").Value = "'1.148"
$ws.Range("E10").Value = "  +1.05%  "
$ws.Range("D11").Value = "'0.07538"
$ws.Range("E11").Value = "  -0.47%  "
$ws.Range("E12").Value = "  -0.08%  "
$ws.Range("D13").Value = "'21.16"
$ws.Range("E13").Value = "  +0.73%  "
$ws.Range("E14").Value = "  +0.86%  "
$ws.Range("D15").Value = "'6.961"
$ws.Range("E15").Value = "  +0.65%  "
$ws.Range("D16").Value = "1.570.65"
$ws.Range("E16").Value = "  +0.50%  "
$ws.Range("E17").Value = "  -0.73%  "
$ws.Range("D18").Value = "'90.47"
$ws.Range("E18").Value = "  +0.76%  "
$ws.Range("D19").Value = "'0.06778"
$ws.Range("E19").Value = "  +0.35%  "
$ws.Range("D21").Value = "'6.338"
$ws.Range("E21").Value = "  +2.24%  "
$ws.Range("E22").Value = "  -1.05%  "
$ws.Range("E23").Value = "  +2.34%  "
$ws.Range("D24").Value = "22.473.19"
$ws.Range("E24").Value = "  +0.48%  "
$ws.Range("D25").Value = "'2.369"
$ws.Range("E25").Value = "  -0.34%  "
$ws.Range("D26").Value = "'2.617"
$ws.Range("E26").Value = "  -3.46%  "
$ws.Range("D27").Value = "'20.02"
$ws.Range("E27").Value = "  -1.00%  "
$ws.Range("D28").Value = "'149.01"
$ws.Range("E28").Value = "  +0.53%  "
$ws.Range("D29").Value = "'5.074"
$ws.Range("E29").Value = "  +0.88%  "
$ws.Range("D31").Value = "1.748.14"
$ws.Range("E31").Value = "  +0.46%  "
$ws.Range("D32").Value = "'1.069"
$ws.Range("E32").Value = "  +7.91%  "
$ws.Range("D33").Value = "'6.191"
$ws.Range("E33").Value = "  +2.34%  "
$ws.Range("D34").Value = "'2.015"
$ws.Range("E34").Value = "  -0.25%  "
$ws.Range("D35").Value = "'9.797"
$ws.Range("E35").Value = "  -2.96%  "
$ws.Range("D36").Value = "'0.08345"
$ws.Range("E36").Value = "  -1.30%  "
$ws.Range("D37").Value = "'1.365"
$ws.Range("E37").Value = "  -4.44%  "
$ws.Range("D38").Value = "'0.02476"
$ws.Range("E38").Value = "  -0.47%  "
$ws.Range("D39").Value = "'0.2300"
$ws.Range("E39").Value = "  +0.39%  "
$ws.Range("E40").Value = "  +1.34%  "
$ws.Range("D41").Value = "'5.428"
$ws.Range("E41").Value = "  +0.22%  "
$ws.Range("D42").Value = "'11.30"
$ws.Range("E42").Value = "  +0.27%  "
$ws.Range("D43").Value = "'0.6208"
$ws.Range("E43").Value = "  -1.65%  "
$ws.Range("E44").Value = "  +1.06%  "
$ws.Range("D46").Value = "'3.807"
$ws.Range("E46").Value = "  +0.26%  "
$ws.Range("D47").Value = "'0.5844"
$ws.Range("E47").Value = "  -1.36%  "
$ws.Range("D48").Value = "'129.29"
$ws.Range("E48").Value = "  +3.65%  "
$ws.Range("E49").Value = "  -0.25%  "
$ws.Range("D50").Value = "'1.228"
$ws.Range("E50").Value = "  -2.77%  "
$ws.Range("D51").Value = "'0.07311"
$ws.Range("E51").Value = "  -0.22%  "
